$wb = $excel.ActiveWorkbook

# --- "Delete" sheet: update the test row that referenced the removed
#     "Terminate after KA Module" parameter/value so it now points at the
#     "Order Take" module and toggles its value from On to Off.
$wsDelete = $wb.Worksheets.Item("Delete")
$wsDelete.Range("B2").Value = "Order Take"
$wsDelete.Range("D2").Value = "Order Take"

# Match E2's border/fill formatting to its neighbour F2 before writing the
# new "Off" value (mirrors the boxed style used by the other toggle cells).
$wsDelete.Range("F2").Copy()
$wsDelete.Range("E2").PasteSpecial(-4122)
$wsDelete.Range("E2").Value = "Off"

# --- Restore selections on each sheet (cursor position captured at save time) ---
$wsCreate = $wb.Worksheets.Item("Create")
$wsCreate.Range("C17").Select()

$wsEdit = $wb.Worksheets.Item("Edit")
$wsEdit.Range("D15").Select()

$wsDelete.Range("D14").Select()

$wsInvalid = $wb.Worksheets.Item("Invalid")
$wsInvalid.Range("C10").Select()

# "Queries" stays the active tab, so select it last.
$wsQueries = $wb.Worksheets.Item("Queries")
$wsQueries.Range("A12").Select()
